# Reorder the "horas" columns on the consolidated-hours report (Verity sheet)
# and tidy up the selection / column width left over from the previous layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Verity")

# --- Row 5 header: swap "Total de horas" / "Banco de horas" and
#     "Horas S.T" / "Horas S.A" (A and E stay put) -------------------------
$ws.Range("B5").Value = "Banco de horas"
$ws.Range("C5").Value = "Total de horas"
$ws.Range("D5").Value = "Horas S.T"
$ws.Range("F5").Value = "Adicional noturno"

# --- Column F is no longer auto-fit; give it a fixed width of 18 ----------
$ws.Columns("F").ColumnWidth = 17.1

# --- Move the active selection up to the header row (A1:F1) ---------------
$ws.Range("A1:F1").Select()
